$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 31.50647122526048
$ws.Range("B3").Value = 18.76444902650253
$ws.Range("B4").Value = 13.54993501083804
$ws.Range("B5").Value = 11.98806451237925
$ws.Range("B6").Value = 11.50634119771969
$ws.Range("B7").Value = 7.015830506684529
$ws.Range("B8").Value = 5.668908520615488
